$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as text-preserving cell values.
# NumberFormat is set to Text ("@") before assignment so that numeric-looking
# strings (e.g. "1.20", "0.999") are preserved exactly instead of being
# auto-coerced into numbers by Excel (which would drop trailing zeros, etc).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.963.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.107.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.61"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.74"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.645.69"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "58.049.36"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.113.90"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.09"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "338.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.66%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.34%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.97"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.13"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.32"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0667"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.155.13"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.46"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.306.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0257"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.70"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.955"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.95%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.17%  "
